$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the now-obsolete "Limit switch" entries (rows 7 & 8, column I) ---
$ws.Range("I7").ClearContents()
$ws.Range("I8").ClearContents()

# --- New relay / PC2 selection-start row ---
$ws.Range("C9").Value = "CN7-35"
$ws.Range("F9").Value = "PC2"
$ws.Range("G9").Value = "input"
$ws.Range("H9").Value = "Start"
$ws.Range("I9").Value = "Selection"

# --- New gas heater start row ---
$ws.Range("C10").Value = "CN7-37"
$ws.Range("F10").Value = "PC3"
$ws.Range("G10").Value = "output"
$ws.Range("H10").Value = "Gas Heater Start"

# --- New I2C temperature sensor rows ---
$ws.Range("C11").Value = "CN10-3"
$ws.Range("D11").Value = "CN5-10"
$ws.Range("E11").Value = "D15"
$ws.Range("F11").Value = "PB8"
$ws.Range("G11").Value = "I2C1_SCL"
$ws.Range("H11").Value = "Temp Sensor"

$ws.Range("C12").Value = "CN10-5"
$ws.Range("D12").Value = "CN5-9"
$ws.Range("E12").Value = "D14"
$ws.Range("F12").Value = "PB9"
$ws.Range("G12").Value = "I2C1_SDA"
$ws.Range("H12").Value = "Temp Sensor"

# --- Column H grew a bit wider to fit the new "Gas Heater Start" text ---
$ws.Columns.Item(8).ColumnWidth = 14.5

# --- The embedded (hidden) Visio drawing shrinks along with column H's resize ---
$shp = $ws.Shapes.Item(1)
$shp.Width = 538.1729

# --- View state: scroll down a bit and leave the selection on the last new cell ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I11").Select()
